$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5501
$ws.Range("J43").Value = 5501
$ws.Range("L43").Value = 5501
$ws.Range("N43").Value = -5639
$ws.Range("H123").Value = 97779.5
$ws.Range("J123").Value = 97779.5
$ws.Range("L123").Value = 97779.5
$ws.Range("N123").Value = -107579.5
$ws.Range("H137").Value = 985390
$ws.Range("I137").Value = 16668492
$ws.Range("J137").Value = 5196.104
$ws.Range("K137").Value = 50005476
$ws.Range("L137").Value = 15588.312
$ws.Range("M137").Value = -50002926
$ws.Range("N137").Value = -20688.312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2996.2222
$ws.Range("I61").Value = 2270.5454
$ws.Range("K61").Value = 2270.5454
$ws.Range("M61").Value = -2058.5454
$ws.Range("H74").Value = 206940.56
$ws.Range("I74").Value = 223247.92
$ws.Range("J74").Value = 3098.5
$ws.Range("K74").Value = 223247.92
$ws.Range("L74").Value = 3098.5
$ws.Range("M74").Value = -222373.92
$ws.Range("N74").Value = -4846.5
$ws.Range("H77").Value = 206940.56
$ws.Range("I77").Value = 223247.92
$ws.Range("J77").Value = 3098.5
$ws.Range("K77").Value = 1116239.6
$ws.Range("L77").Value = 15492.5
$ws.Range("M77").Value = -1111871.6
$ws.Range("N77").Value = -24228.5
$ws.Range("H132").Value = 5466959.5
$ws.Range("I132").Value = 2200.0408
$ws.Range("J132").Value = 27781396
$ws.Range("K132").Value = 6600.1224
$ws.Range("L132").Value = 83344188
$ws.Range("M132").Value = -4070.1224
$ws.Range("N132").Value = -83349248
$ws.Range("H136").Value = 2996.2222
$ws.Range("I136").Value = 2270.5454
$ws.Range("K136").Value = 6811.6362
$ws.Range("M136").Value = -4261.6362
$ws.Range("H141").Value = 93250
$ws.Range("J141").Value = 93250
$ws.Range("L141").Value = 93250
$ws.Range("N141").Value = -103610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2788.3044
$ws.Range("I134").Value = 2449.258
$ws.Range("J134").Value = 3489
$ws.Range("K134").Value = 7347.773999999999
$ws.Range("L134").Value = 10467
$ws.Range("M134").Value = -4812.773999999999
$ws.Range("N134").Value = -15537

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5525.1665
$ws.Range("I31").Value = 3372.7334
$ws.Range("K31").Value = 3372.7334
$ws.Range("M31").Value = -3077.7334
$ws.Range("H34").Value = 5525.1665
$ws.Range("I34").Value = 3372.7334
$ws.Range("K34").Value = 3372.7334
$ws.Range("M34").Value = -3170.7334
$ws.Range("H58").Value = 2066.75
$ws.Range("I58").Value = 1669.5
$ws.Range("J58").Value = 2622.9
$ws.Range("K58").Value = 1669.5
$ws.Range("L58").Value = 2622.9
$ws.Range("M58").Value = -1466.5
$ws.Range("N58").Value = -3028.9
$ws.Range("H74").Value = 79314
$ws.Range("J74").Value = 79314
$ws.Range("L74").Value = 79314
$ws.Range("N74").Value = -81062
$ws.Range("H77").Value = 79314
$ws.Range("J77").Value = 79314
$ws.Range("L77").Value = 237942
$ws.Range("N77").Value = -246678
$ws.Range("H132").Value = 12348858
$ws.Range("I132").Value = 2197.6875
$ws.Range("K132").Value = 6593.0625
$ws.Range("M132").Value = -4063.0625
$ws.Range("H136").Value = 2066.75
$ws.Range("I136").Value = 1669.5
$ws.Range("J136").Value = 2622.9
$ws.Range("K136").Value = 5008.5
$ws.Range("L136").Value = 7868.700000000001
$ws.Range("M136").Value = -2458.5
$ws.Range("N136").Value = -12968.7
$ws.Range("H141").Value = 447824
$ws.Range("J141").Value = 447824
$ws.Range("L141").Value = 447824
$ws.Range("N141").Value = -458184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3972.8333
$ws.Range("J81").Value = 4586
$ws.Range("L81").Value = 13758
$ws.Range("N81").Value = -16004
$ws.Range("H84").Value = 3972.8333
$ws.Range("J84").Value = 4586
$ws.Range("L84").Value = 41274
$ws.Range("N84").Value = -52506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2493.3584
$ws.Range("I122").Value = 1924.0731
$ws.Range("J122").Value = 4438.4165
$ws.Range("K122").Value = 5772.219300000001
$ws.Range("L122").Value = 13315.2495
$ws.Range("M122").Value = -3322.219300000001
$ws.Range("N122").Value = -18215.2495
$ws.Range("H132").Value = 2542.7576
$ws.Range("J132").Value = 4819.25
$ws.Range("L132").Value = 14457.75
$ws.Range("N132").Value = -19517.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 363
$ws.Range("J22").Value = 454.75
$ws.Range("K22").Value = 363
$ws.Range("L22").Value = 454.75
$ws.Range("M22").Value = -68
$ws.Range("N22").Value = -1044.75
$ws.Range("I27").Value = 363
$ws.Range("J27").Value = 454.75
$ws.Range("K27").Value = 363
$ws.Range("L27").Value = 454.75
$ws.Range("M27").Value = -256
$ws.Range("N27").Value = -668.75
$ws.Range("H55").Value = 516.75
$ws.Range("J55").Value = 1440
$ws.Range("L55").Value = 1440
$ws.Range("N55").Value = -1786
$ws.Range("H132").Value = 4880.9287
$ws.Range("J132").Value = 7336.091
$ws.Range("L132").Value = 22008.273
$ws.Range("N132").Value = -27068.273
$ws.Range("H136").Value = 2145.577
$ws.Range("I136").Value = 1945.9535
$ws.Range("K136").Value = 5837.860500000001
$ws.Range("M136").Value = -3287.860500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 16490
$ws.Range("I58").Value = 16490
$ws.Range("K58").Value = 16490
$ws.Range("M58").Value = -16182
$ws.Range("H113").Value = 1236.6552
$ws.Range("I113").Value = 1028.0555
$ws.Range("K113").Value = 3084.1665
$ws.Range("M113").Value = -914.1664999999998
$ws.Range("H122").Value = 25002802
$ws.Range("I122").Value = 3167.1667
$ws.Range("K122").Value = 9501.500100000001
$ws.Range("M122").Value = -7051.500100000001
$ws.Range("H132").Value = 3651.5652
$ws.Range("I132").Value = 4318.8
$ws.Range("J132").Value = 3138.3076
$ws.Range("K132").Value = 12956.4
$ws.Range("L132").Value = 9414.9228
$ws.Range("M132").Value = -10426.4
$ws.Range("N132").Value = -14474.9228
$ws.Range("H136").Value = 5700.7886
$ws.Range("I136").Value = 5318.3076
$ws.Range("J136").Value = 6848.231
$ws.Range("K136").Value = 15954.9228
$ws.Range("L136").Value = 20544.693
$ws.Range("M136").Value = -13404.9228
$ws.Range("N136").Value = -25644.693
